$data = @(
  @("ECs","Ntn3","Neo1","ECs",3,1,0.248739,0.7462169999999999,0.09500454514204523,0.09500454514204522,3,1,5.872120333333332,17.616361,0.07819433676692768,0.07819433676692769,1.460625339593,13.145628056337,0.007428817397225868,0.007428817397225868),
  @("ECs","Ntn3","Neo1","FAPs",3,1,0.248739,0.7462169999999999,0.09500454514204523,0.09500454514204522,3,1,59.62659933333333,178.879798,0.7939998031155241,0.7939998031155242,14.831460691574,133.483146224166,0.07543359013786384,0.07543359013786384),
  @("ECs","Ntn3","Neo1","M2",3,1,0.248739,0.7462169999999999,0.09500454514204523,0.09500454514204522,3,1,0.226786,0.680358,0.003019928041555994,0.003019928041555994,0.05641052285399999,0.5076947056859999,0.0002869068899497347,0.0002869068899497346),
  @("ECs","Ntn3","Neo1","sCs",3,1,0.248739,0.7462169999999999,0.09500454514204523,0.09500454514204522,3,1,9.370985666666668,28.112957,0.1247859320759922,0.1247859320759922,2.330929603741,20.978366433669,0.01185523071700579,0.01185523071700579),
  @("FAPs","Ntn3","Neo1","ECs",3,1,1.982221333333333,5.946664,0.7570989516890869,0.7570989516890868,3,1,5.872120333333332,17.616361,0.07819433676692768,0.07819433676692769,11.63984219663378,104.758579769704,0.05920085039426437,0.05920085039426437),
  @("FAPs","Ntn3","Neo1","FAPs",3,1,1.982221333333333,5.946664,0.7570989516890869,0.7570989516890868,3,1,59.62659933333333,178.879798,0.7939998031155241,0.7939998031155242,118.1931172326524,1063.738055093872,0.6011364185801047,0.6011364185801047),
  @("FAPs","Ntn3","Neo1","M2",3,1,1.982221333333333,5.946664,0.7570989516890869,0.7570989516890868,3,1,0.226786,0.680358,0.003019928041555994,0.003019928041555994,0.4495400473013334,4.045860425712,0.00228638435443852,0.00228638435443852),
  @("FAPs","Ntn3","Neo1","sCs",3,1,1.982221333333333,5.946664,0.7570989516890869,0.7570989516890868,3,1,9.370985666666668,28.112957,0.1247859320759922,0.1247859320759922,18.57536770282756,167.178309325448,0.09447529836027928,0.09447529836027926),
  @("M2","Ntn3","Neo1","ECs",1,0.3333333333333333,0.07456866666666666,0.223706,0.02848110774151,0.02848110774151,3,1,5.872120333333332,17.616361,0.07819433676692768,0.07819433676692769,0.4378761837628888,3.940885653865999,0.002227061330234784,0.002227061330234785),
  @("M2","Ntn3","Neo1","FAPs",1,0.3333333333333333,0.07456866666666666,0.223706,0.02848110774151,0.02848110774151,3,1,59.62659933333333,178.879798,0.7939998031155241,0.7939998031155242,4.446276010154222,40.01648409138799,0.02261399393927097,0.02261399393927098),
  @("M2","Ntn3","Neo1","M2",1,0.3333333333333333,0.07456866666666666,0.223706,0.02848110774151,0.02848110774151,3,1,0.226786,0.680358,0.003019928041555994,0.003019928041555994,0.01691112963866667,0.152200166748,0.00008601089592316358,0.00008601089592316358),
  @("M2","Ntn3","Neo1","sCs",1,0.3333333333333333,0.07456866666666666,0.223706,0.02848110774151,0.02848110774151,3,1,9.370985666666668,28.112957,0.1247859320759922,0.1247859320759922,0.6987819065157778,6.289037158642,0.003554041576081082,0.003554041576081082),
  @("sCs","Ntn3","Neo1","ECs",3,1,0.312651,0.9379529999999999,0.1194153954273579,0.1194153954273579,3,1,5.872120333333332,17.616361,0.07819433676692768,0.07819433676692769,1.835924294336999,16.523318649033,0.009337607645202662,0.009337607645202662),
  @("sCs","Ntn3","Neo1","FAPs",3,1,0.312651,0.9379529999999999,0.1194153954273579,0.1194153954273579,3,1,59.62659933333333,178.879798,0.7939998031155241,0.7939998031155242,18.642315908166,167.780843173494,0.09481580045828465,0.09481580045828465),
  @("sCs","Ntn3","Neo1","M2",3,1,0.312651,0.9379529999999999,0.1194153954273579,0.1194153954273579,3,1,0.226786,0.680358,0.003019928041555994,0.003019928041555994,0.07090486968599999,0.638143827174,0.0003606259012445757,0.0003606259012445756),
  @("sCs","Ntn3","Neo1","sCs",3,1,0.312651,0.9379529999999999,0.1194153954273579,0.1194153954273579,3,1,9.370985666666668,28.112957,0.1247859320759922,0.1247859320759922,2.929848039669,26.368632357021,0.01490136142262603,0.01490136142262603)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $i + 2
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $row[$j]
  }
}
